$wb = $excel.ActiveWorkbook

# --- Arkusz2 (sheet3.xml): add helper ratio cells Q3:S3, highlighted yellow ---
# (done first so that the workbook ends up with Sheet1 as the active/selected
# sheet, matching the final saved view state)
$ws3 = $wb.Worksheets.Item("Arkusz2")
$ws3.Range("R3").Value = 0
$ws3.Range("S3").Value = 2
$ws3.Range("R3:S3").Interior.Color = 65535
$ws3.Range("Q3").Formula = "=R3/S3"
$ws3.Range("R4").Select()

# --- Sheet1: add two new rows (59 and 60) to the task list ---
$ws1 = $wb.Worksheets.Item("Sheet1")

# Row 59
$ws1.Range("A59").Value = 58
$ws1.Range("A59").VerticalAlignment = -4160
$ws1.Range("B59").Value = "wyliczenia"
$ws1.Range("B59").VerticalAlignment = -4160
$ws1.Range("C59").Value = "Stworzenie tabel tymczasowych dla wyliczeń. Ilość tabel niezbędna do zapisania wszystkich zestawień. Tabele te są tymczasowe i mogą być czyszczone zawsze podczas otwierania kolejnego miesiąca."
$ws1.Range("C59").WrapText = $true
$ws1.Range("D59").Value = 0
$ws1.Range("D59").VerticalAlignment = -4108
$ws1.Rows.Item(59).RowHeight = 30

# Row 60
$ws1.Range("A60").Value = 59
$ws1.Range("A60").VerticalAlignment = -4160
$ws1.Range("B60").Value = "Wyliczenia mag"
$ws1.Range("B60").VerticalAlignment = -4160
$ws1.Range("C60").Value = "dla WYDANIA nieobecności są indywidualne i wymagają oddzielnej procedury SQL - należy przerobić kod by ujednolicić procedurę. "
$ws1.Range("C60").WrapText = $true
$ws1.Range("D60").Value = 0
$ws1.Range("D60").VerticalAlignment = -4108
$ws1.Rows.Item(60).RowHeight = 30

# Cursor / view position ends on Sheet1, matching the authored selection
$ws1.Range("B61").Select()

$excel.Calculate()
